# Updates the crypto price table (Sheet1) to the new snapshot described by
# the diff: a handful of price (column D) refreshes, two coins ("LEO" and
# "One") moving to a new rank, and every row between them shifting down by
# one, plus their Volume(1h) text (column E) and rank-prefixed strings
# being recomputed accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the value to be stored as text, preserving formats like
    # trailing zeros ("243.38") instead of being auto-converted to a number.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "243.38"
Set-TextValue $ws.Range("D3") "23.05"
$ws.Range("B4").Value = "LEO"
$ws.Range("C4").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D4") "3.637"
$ws.Range("E4").Value = "3LEOLEO"
$ws.Range("B5").Value = "HuobiToken"
$ws.Range("C5").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D5") "5.399"
$ws.Range("E5").Value = "4HuobiTokenHT"
$ws.Range("B6").Value = "Cronos"
$ws.Range("C6").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D6") "0.05935"
$ws.Range("E6").Value = "5CronosCRO"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D7") "3.400"
$ws.Range("E7").Value = "6GateTokenGT"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D8") "0.8087"
$ws.Range("E8").Value = "7MXTokenMX"
$ws.Range("B9").Value = "FTXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D9") "0.9135"
$ws.Range("E9").Value = "8FTXTokenFTT"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D10") "0.01114"
$ws.Range("E10").Value = "9OneONEBestin24h"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D11") "0.1415"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D12") "0.07422"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D13") "0.03334"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D14") "0.03066"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D15") "0.09332"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws.Range("D16") "3.949"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D17") "0.001582"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws.Range("D18") "0.04803"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D19") "0.005436"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue $ws.Range("D20") "0.004426"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue $ws.Range("D21") "0.0009861"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue $ws.Range("D22") "0.00007809"
$ws.Range("E22").Value = "21NitroExNTX"
Set-TextValue $ws.Range("D23") "6.426"
Set-TextValue $ws.Range("D26") "0.1348"
Set-TextValue $ws.Range("D40") "0.03874"
Set-TextValue $ws.Range("D41") "0.006218"
$ws.Range("E41").Value = "40KickTokenKICK"
Set-TextValue $ws.Range("D42") "0.1068"
Set-TextValue $ws.Range("D43") "0.002903"
Set-TextValue $ws.Range("D44") "0.006666"
Set-TextValue $ws.Range("D45") "0.00005177"
Set-TextValue $ws.Range("D47") "0.0005805"
Set-TextValue $ws.Range("D48") "0.8510"
Set-TextValue $ws.Range("D49") "0.002271"
Set-TextValue $ws.Range("D50") "0.00002102"
